# Apply the edit described by the commit:
# "Add check if the property defined in the excel file is supported by the
#  configuration object"
#
# Concretely (derived from the OOXML diff):
#  - On the "DataCombined" sheet, fill in the missing `group` value
#    ("Aciclovir PVB") for the two AciclovirPVB rows (G2, G3).
#  - On the "plotConfiguration" sheet, remove the stray/unsupported test
#    values that had been typed into column N (N1="blabla", N2="1;2;3") -
#    these were outside the real table and not a supported property.

$wb = $excel.ActiveWorkbook

$dataCombined = $wb.Worksheets.Item("DataCombined")
$dataCombined.Activate()
$dataCombined.Range("G2").Value = "Aciclovir PVB"
$dataCombined.Range("G3").Value = "Aciclovir PVB"
$dataCombined.Range("G3").Select()

$plotConfiguration = $wb.Worksheets.Item("plotConfiguration")
$plotConfiguration.Activate()
$plotConfiguration.Range("N1").ClearContents()
$plotConfiguration.Range("N2").ClearContents()

$dataCombined.Activate()

$wb.Save()
